$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "가치"
$ws.Range("A3").Value = "감사"
$ws.Range("A4").Value = "걱정"
$ws.Range("A5").Value = "건강"
$ws.Range("A6").Value = "결정"
$ws.Range("A7").Value = "경험"
$ws.Range("A8").Value = "계획"
$ws.Range("A9").Value = "고민"
$ws.Range("A10").Value = "고통"
$ws.Range("A11").Value = "공부"
$ws.Range("A12").Value = "관심"
$ws.Range("A13").Value = "교육"
$ws.Range("A14").Value = "기대"
$ws.Range("A15").Value = "기분"
$ws.Range("A16").Value = "기쁨"
$ws.Range("A17").Value = "기억"
$ws.Range("A18").Value = "기적"
$ws.Range("A19").Value = "기회"
$ws.Range("A20").Value = "나가다"
$ws.Range("A21").Value = "나누다"
$ws.Range("A22").Value = "나이"
$ws.Range("A23").Value = "나타나다"
$ws.Range("A24").Value = "남기다"
$ws.Range("A25").Value = "남다"
$ws.Range("A26").Value = "남자"
$ws.Range("A27").Value = "낮다"
$ws.Range("A28").Value = "내려놓다"
$ws.Range("A29").Value = "내리다"
$ws.Range("A30").Value = "내일"
$ws.Range("A31").Value = "너"
$ws.Range("A32").Value = "넘다"
$ws.Range("A33").Value = "넘어오다"
$ws.Range("A34").Value = "노력"
$ws.Range("A35").Value = "놀다"
$ws.Range("A36").Value = "놀라다"
$ws.Range("A37").Value = "놓치다"
$ws.Range("A38").Value = "누구"
$ws.Range("A39").Value = "늘어나다"
$ws.Range("A40").Value = "다리"
$ws.Range("A41").Value = "다음"
$ws.Range("A42").Value = "닫다"
$ws.Range("A43").Value = "대화"
$ws.Range("A44").Value = "덕분에"
$ws.Range("A45").Value = "도움"
$ws.Range("A46").Value = "독서"
$ws.Range("A47").Value = "아끼다"
$ws.Range("A48").Value = "아름답다"
$ws.Range("A49").Value = "안녕하세요"
$ws.Range("A50").Value = "안심"
$ws.Range("A51").Value = "안타깝다"
$ws.Range("A52").Value = "알다"
$ws.Range("A53").Value = "알리다(알려주다)"
$ws.Range("A54").Value = "알아서하다"
$ws.Range("A55").Value = "알아차리다"
$ws.Range("A56").Value = "약속"
$ws.Range("A57").Value = "어떻게"
$ws.Range("A58").Value = "어렵다"
$ws.Range("A59").Value = "어린이"
$ws.Range("A60").Value = "어지럽다"
$ws.Range("A61").Value = "오늘"
$ws.Range("A62").Value = "온도"
$ws.Range("A63").Value = "왜"
$ws.Range("A64").Value = "웃다"
$ws.Range("A65").Value = "원하다"
$ws.Range("A66").Value = "이름"
$ws.Range("A67").Value = "이해"
$ws.Range("A68").Value = "일부러"
$ws.Range("A69").Value = "잃어버리다"
$ws.Range("A70").Value = "입원"
$ws.Range("A71").Value = "잊다"
